$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.690.49"
$ws.Range("E2").Value = "  +0.82%  "

$ws.Range("D3").Value = "2.277.39"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'503.91"
$ws.Range("E5").Value = "  +1.70%  "

$ws.Range("D6").Value = "'128.54"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "'0.528"
$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("D9").Value = "2.291.35"
$ws.Range("E9").Value = "  +0.37%  "

$ws.Range("D10").Value = "'0.0963"
$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").Value = "'0.345"
$ws.Range("E12").Value = "  +4.55%  "

$ws.Range("D13").Value = "'4.91"
$ws.Range("E13").Value = "  +4.27%  "

$ws.Range("D14").Value = "'23.33"
$ws.Range("E14").Value = "  +5.04%  "

$ws.Range("D15").Value = "2.682.58"
$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("D16").Value = "54.739.56"
$ws.Range("E16").Value = "  +1.02%  "

$ws.Range("D17").Value = "'0.0000130"
$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("D18").Value = "2.288.50"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("D19").Value = "'10.29"
$ws.Range("E19").Value = "  +1.96%  "

$ws.Range("D20").Value = "'4.13"
$ws.Range("E20").Value = "  +0.60%  "

$ws.Range("D21").Value = "'306.80"
$ws.Range("E21").Value = "  +1.09%  "

$ws.Range("E22").Value = "  -0.45%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'60.05"
$ws.Range("E24").Value = "  -2.79%  "

$ws.Range("D25").Value = "'0.995"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("D27").Value = "'7.44"
$ws.Range("E27").Value = "  +3.04%  "

$ws.Range("D28").Value = "'170.83"
$ws.Range("E28").Value = "  +0.22%  "

$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'6.06"
$ws.Range("E29").Value = "  +3.17%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0702"
$ws.Range("E30").Value = "  +3.09%  "

$ws.Range("D31").Value = "'1.62"
$ws.Range("E31").Value = "  +0.92%  "

$ws.Range("E32").Value = "  +2.14%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Value = "'17.93"
$ws.Range("E34").Value = "  +1.31%  "

$ws.Range("D35").Value = "'0.994"
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("D36").Value = "'0.917"
$ws.Range("E36").Value = "  +2.08%  "

$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("D38").Value = "'3.77"
$ws.Range("E38").Value = "  +2.08%  "

$ws.Range("D39").Value = "'36.32"
$ws.Range("E39").Value = "  +1.39%  "

$ws.Range("D40").Value = "'0.375"
$ws.Range("E40").Value = "  +0.49%  "

$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").Value = "'5.04"
$ws.Range("E42").Value = "  +5.51%  "

$ws.Range("D43").Value = "'3.39"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").Value = "'125.95"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "'0.0497"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("D46").Value = "'248.20"
$ws.Range("E46").Value = "  +4.28%  "

$ws.Range("D47").Value = "'0.0901"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("D48").Value = "'0.548"
$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("D49").Value = "'0.374"
$ws.Range("E49").Value = "  +0.80%  "

$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").Value = "'10.81"
